# Glossary update: add PT/EN/ES language rows to the table and extend print setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New glossary rows (Acronym / Term)
# Shared-string insertion order must match: PT, EN, ES, Portuguese, English, Spanish
$ws.Range("A18").Value = "PT"
$ws.Range("A19").Value = "EN"
$ws.Range("A20").Value = "ES"

$ws.Range("B18").Value = "Portuguese"
$ws.Range("B19").Value = "English"
$ws.Range("B20").Value = "Spanish"

# Apply the same style used by the rest of the Acronym/TEA columns (bold, vertically centered)
$ws.Range("A18:B20").Font.Bold = $true
$ws.Range("A18:B20").VerticalAlignment = -4108

# Resize the Excel Table (ListObject) to include the new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:C20"))

# Update the visible window / selection to reflect scrolling to the new rows
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C25").Select()

# Configure page setup for printing
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = 300
